# "Generate Report for Archive"
#
# 1) Shared string update: "Ready for handoff" -> "In Translation"
#    (affects every cell across all three sheets that shows this status:
#     Overview!E2:E3/F2:F3, zh-cn!C2:C3, de-de!C2:C3)
# 2) Narrow the "status" columns (E/F on Overview, C on zh-cn / de-de)
#    from ~17.22 chars down to ~13.41 chars.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item(1)
$zhcn     = $wb.Worksheets.Item(2)
$dede     = $wb.Worksheets.Item(3)

# --- 1) Replace the status text everywhere it appears ---------------------
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"

$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"

$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"

# --- 2) Narrow the columns -------------------------------------------------
# ColumnWidth is in characters and gets quantized to the engine's internal
# pixel grid (1/6-character steps), so 12.5 is the closest input that lands
# on the nearest achievable width to the target 13.4101848602295.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
